$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 25; this shifts the existing rows 25-35 down to 26-36
# and copies formatting (incl. the date style on column D) down from the row
# that used to be at 25.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Cells.Item(25, 1).Value = 7
$ws.Cells.Item(25, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(25, 3).Value = "Ñuble"
$ws.Cells.Item(25, 4).Value = (Get-Date -Year 2022 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(25, 5).Value = 16
$ws.Cells.Item(25, 6).Value = 100112043
$ws.Cells.Item(25, 7).Value = "Pepino dulce"
$ws.Cells.Item(25, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 100
$ws.Cells.Item(25, 11).Value = 17000
$ws.Cells.Item(25, 12).Value = 18000
$ws.Cells.Item(25, 13).Value = 17500
$ws.Cells.Item(25, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(25, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(25, 16).Value = 972
$ws.Cells.Item(25, 17).Value = 18
$ws.Cells.Item(25, 18).Value = "Hortaliza"
